$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell E8 previously displayed "Good Morning"; update it to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Selection moved to E8 on the active sheet
$ws.Range("E8").Select()
